# Weekly Hortaliza / Albahaca price log update.
# A new weekly record is inserted as row 77 (pushing the former rows 77 and
# 78 down to 78 and 79 respectively), growing the sheet from 78 to 79 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the old row 77; this shifts old rows 77/78 -> 78/79
# and automatically extends the used range / dimension.
$ws.Rows.Item(77).Insert()

# Populate the newly inserted row 77 with the new weekly record.
$ws.Cells.Item(77, 1).Value = 8
$ws.Cells.Item(77, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(77, 3).Value = "Coquimbo"
$ws.Cells.Item(77, 4).Value = 44595
$ws.Cells.Item(77, 5).Value = 4
$ws.Cells.Item(77, 6).Value = 100112052
$ws.Cells.Item(77, 7).Value = "Albahaca"
$ws.Cells.Item(77, 8).Value = "Sin especificar"
$ws.Cells.Item(77, 9).Value = "Primera"
$ws.Cells.Item(77, 10).Value = 800
$ws.Cells.Item(77, 11).Value = 3500
$ws.Cells.Item(77, 12).Value = 4000
$ws.Cells.Item(77, 13).Value = 3750
$ws.Cells.Item(77, 14).Value = "`$/docena de matas"
$ws.Cells.Item(77, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(77, 16).Value = 625
$ws.Cells.Item(77, 17).Value = 6
$ws.Cells.Item(77, 18).Value = "Hortaliza"
